# 2006 Monthly Time Charts
# Adds a new "monthly-time-chart" worksheet after the existing
# "drug-arrests-by-age-and-type-of" sheet, with a Month header row and a
# single "Number of Drug Arrests" data row (Jan-Dec), formatted with a
# thousands-separator number format, right aligned.

$wb = $excel.ActiveWorkbook
$existing = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add($null, $existing)
$ws.Name = "monthly-time-chart"

# --- Header row: Month, January .. December -----------------------------
$months = @("Month","January","February","March","April","May","June","July","August","September","October","November","December")
for ($i = 0; $i -lt $months.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $months[$i]
}

# --- Data row: Number of Drug Arrests ------------------------------------
$ws.Range("A2").Value = "Number of Drug Arrests"

$counts = @(2714, 2788, 3119, 2581, 2900, 2770, 2664, 2802, 2551, 2559, 2301, 2591)
for ($i = 0; $i -lt $counts.Length; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $counts[$i]
}

# Thousands-separator number format, right aligned, on the monthly values
$dataRange = $ws.Range("B2:M2")
$dataRange.NumberFormat = "#,###"
$dataRange.HorizontalAlignment = -4152   # xlRight

# Column A a bit wider, to match the label column on the other sheet
$ws.Columns.Item(1).ColumnWidth = 31.43

$existing.Activate()

$wb.Save()
